$d = $word.ActiveDocument

# Locate the "Logo Ametrine Archive" picture (the second inline picture in the
# document) and anchor right after the paragraph that holds it. The new
# "PERENCANAAN WEBSITE" section header block gets inserted there, before the
# existing blank paragraphs / page break that lead into "JADWAL PELAKSANAAN
# PROYEK".
$shape = $d.InlineShapes.Item(2)
$imgParagraph = $shape.Range.Paragraphs.Item(1)
$insertPoint = $d.Range($imgParagraph.Range.End, $imgParagraph.Range.End)

# Four new paragraphs: blank line, a page break, the centered bold title
# "PERENCANAAN WEBSITE", and a trailing blank centered line. Supplying them
# as one multi-paragraph OOXML fragment (rather than one <w:p> at a time)
# makes Word splice in genuinely new paragraphs ahead of the existing
# "0F895E4F" blank paragraph instead of merging text into it.
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:spacing w:line="360" w:lineRule="auto"/><w:jc w:val="both"/><w:rPr><w:rFonts w:asciiTheme="majorBidi" w:hAnsiTheme="majorBidi" w:cstheme="majorBidi"/><w:b/><w:bCs/></w:rPr></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:rFonts w:asciiTheme="majorBidi" w:hAnsiTheme="majorBidi" w:cstheme="majorBidi"/><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:asciiTheme="majorBidi" w:hAnsiTheme="majorBidi" w:cstheme="majorBidi"/><w:b/><w:bCs/></w:rPr><w:br w:type="page"/></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:spacing w:line="360" w:lineRule="auto"/><w:jc w:val="center"/><w:rPr><w:rFonts w:asciiTheme="majorBidi" w:hAnsiTheme="majorBidi" w:cstheme="majorBidi"/><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:asciiTheme="majorBidi" w:hAnsiTheme="majorBidi" w:cstheme="majorBidi"/><w:b/><w:bCs/></w:rPr><w:lastRenderedPageBreak/><w:t>PERENCANAAN WEBSITE</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:spacing w:line="360" w:lineRule="auto"/><w:jc w:val="center"/><w:rPr><w:rFonts w:asciiTheme="majorBidi" w:hAnsiTheme="majorBidi" w:cstheme="majorBidi"/><w:b/><w:bCs/></w:rPr></w:pPr></w:p>'

$insertPoint.InsertXML($xml) | Out-Null
